{"js": "// The document had two menu passages that were fully upper-case and were\n// changed to lower-case (everything else in the observable diff is just\n// Word's grammar-checker splitting runs / inserting <w:proofErr/> markers\n// around the same, unchanged text \u2014 i.e. no visible content change).\n//\n//   1) \"FIDEOS DE TRIGO, HUEVO COCIDO, WAKAME, SOYA Y PROTE\u00cdNA\"\n//        -> \"fideos de trigo, huevo cocido, wakame, soya y prote\u00edna\"\n//   2) \"TROCITOS DE POLLO APANADOS CON PAPAS RIZADAS\"\n//        -> \"trocitos de pollo apanados con papas rizadas\"\n//      \"Y JUGO HIT EN CAJA\"  (rest of that paragraph, the dot leaders and\n//      price, is untouched since it has no letters to case-fold)\n//        -> \"y jugo hit en caja\"\n\nasync function lowerFirstMatch(searchText) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    const found = results.items[0];\n    found.insertText(found.text.toLowerCase(), Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\nawait lowerFirstMatch(\"FIDEOS DE TRIGO, HUEVO COCIDO, WAKAME, SOYA Y PROTE\u00cdNA\");\nawait lowerFirstMatch(\"TROCITOS DE POLLO APANADOS CON PAPAS RIZADAS\");\nawait lowerFirstMatch(\"Y JUGO HIT EN CAJA\");\n", "ps1": "# The document had two menu passages that were fully upper-case and were\n# changed to lower-case (everything else in the observable diff is just\n# Word's grammar-checker splitting runs / inserting proofErr markers around\n# the same, unchanged text -- i.e. no visible content change).\n#\n#   1) \"FIDEOS DE TRIGO, HUEVO COCIDO, WAKAME, SOYA Y PROTE\u00cdNA\"\n#        -> \"fideos de trigo, huevo cocido, wakame, soya y prote\u00edna\"\n#   2) \"TROCITOS DE POLLO APANADOS CON PAPAS RIZADAS\"\n#        -> \"trocitos de pollo apanados con papas rizadas\"\n#      \"Y JUGO HIT EN CAJA\" (rest of that paragraph, the dot leaders and\n#      price, is untouched since it has no letters to case-fold)\n#        -> \"y jugo hit en caja\"\n\n$doc = $word.ActiveDocument\n\nfunction Lower-FirstMatch($searchText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $searchText.ToLower()\n    $find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\nLower-FirstMatch \"FIDEOS DE TRIGO, HUEVO COCIDO, WAKAME, SOYA Y PROTE\u00cdNA\"\nLower-FirstMatch \"TROCITOS DE POLLO APANADOS CON PAPAS RIZADAS\"\nLower-FirstMatch \"Y JUGO HIT EN CAJA\"\n"}
